$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New review rows to append (appid, keyword, email1, email2, recovery/time, review, blue)
$rows = @(
    @{ C = "orenatias858@gmail.com"; D = "dan624655@gmail.com"; F = "awesome topics about bitcoin" },
    @{ C = "rabuhav25@gmail.com "; D = "itaisenior@gmail.com"; F = "all what I needed to know in one place" },
    @{ C = "innaplutov1@gmail.com"; D = "rabuhav25@gmail.com "; F = "App is working as expected after update.." }
)

$startRow = 17
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Value = "com.hamxa.shaynachim"
    $cellA.Font.Name = "Mangal"

    $cellB = $ws.Cells.Item($r, 2)
    $cellB.Value = "bitcoin"
    $cellB.Font.Name = "Arial"

    $ws.Cells.Item($r, 3).Value = $data.C
    $ws.Cells.Item($r, 4).Value = $data.D

    $cellE = $ws.Cells.Item($r, 5)
    $cellE.Value = "27/5/2019 15:59"
    $cellE.Font.Name = "Arial"

    $cellF = $ws.Cells.Item($r, 6)
    $cellF.Value = $data.F
    $cellF.Font.Name = "Arial"

    $cellG = $ws.Cells.Item($r, 7)
    $cellG.Value = "yes"
    $cellG.Font.Name = "Arial"
}

# Add hyperlink for the new itaisenior@gmail.com address in D18, then restore
# the cell's original (pre-hyperlink) look, since Excel auto-applies its
# built-in "Hyperlink" style when a hyperlink is inserted.
[void]$ws.Hyperlinks.Add($ws.Range("D18"), "mailto:itaisenior@gmail.com", "", "", "itaisenior@gmail.com")
$wb.Styles.Item("Hyperlink").Delete()
$linkFont = $ws.Range("D18").Font
$linkFont.Underline = -4142
$linkFont.Color = 0
$linkFont.Name = "Calibri"
$linkFont.Size = 11

# Move the active selection to G20
[void]$ws.Range("G20").Select()
